# Apply cryptocurrency price/volume updates to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell holding the untouched default (style 0) formatting —
# used to strip the quote-prefix/text style Excel applies when we force
# numeric-looking strings (e.g. "1.000") to stay text instead of becoming numbers.
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "30.259.54"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.863.94"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = $defaultStyle
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.45"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4683"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = "  +0.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2864"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06545"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.19"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  +11.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07895"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.82"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").Value = "1.868.34"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.186"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6803"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  +1.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "278.10"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("D17").Value = "30.258.07"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("E18").Value = "  +7.69%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007340"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.382"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  -2.54%  "
$ws.Range("D22").Value = "2.110.92"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.193"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.45"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.280"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.07"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.944"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.385"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  +3.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09814"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  +2.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.383"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.481"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.069"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04750"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.142"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  +4.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7055"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.705"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.622"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "  +4.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "76.40"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  +3.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.292"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.956"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  +1.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8509"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4181"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.13"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.219"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "949.18"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "  -4.67%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.243"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.26"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05637"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  -0.02%  "

Write-Output "Applied 99 cell updates"
